# feat: add 2022-Q4 data
#
# Before: sheets = [总计, 2022-Q3]
# After:  sheets = [总计, 2022-Q4, 2022-Q3]
#   - 总计: new row inserted summarising the 2022-Q4 quarter (pushes the
#     existing 2022-Q3 summary row down one row).
#   - 2022-Q4: brand-new worksheet (fund holdings table for the new quarter),
#     inserted right before the existing 2022-Q3 sheet.
#   - 2022-Q3: the pre-existing worksheet, left as-is (just shifted right).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" (summary) sheet: push the current data row down and write the
#    new 2022-Q4 summary into the row it vacates.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Copy row 2 (value + formatting) down into row 3 first, so the old
# 2022-Q3 summary keeps its original styling once we overwrite row 2.
$summary.Range("A2:D2").Copy($summary.Range("A3:D3")) | Out-Null

# Row 3 now holds a duplicate of the old 2022-Q3 summary; fix up A3 (the
# running index) and re-assert the rest explicitly.
$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(3, 2).Value = "2022-Q3"
$summary.Cells.Item(3, 3).Value = 2
$summary.Cells.Item(3, 4).Value = 0.52

# Row 2 becomes the new 2022-Q4 summary (A2 stays 0).
$summary.Cells.Item(2, 2).Value = "2022-Q4"
$summary.Cells.Item(2, 3).Value = 5
$summary.Cells.Item(2, 4).Value = 0.73

# ---------------------------------------------------------------------
# 2) Insert the new "2022-Q4" worksheet right before "2022-Q3".
# ---------------------------------------------------------------------
$oldQ3 = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($oldQ3)
$q4.Name = "2022-Q4"

# Helper: write a value that must be preserved as TEXT even though it
# looks numeric (e.g. "4.71"), matching the source data's inlineStr type.
function Set-TextCell($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Header row - match the bold/centered/bordered style used by the other
# header rows in this workbook (copy format from the "总计" header cell).
$summary.Range("B1").Copy($q4.Range("B1:H1")) | Out-Null
$q4.Cells.Item(1, 2).Value = "基金代码"
$q4.Cells.Item(1, 3).Value = "基金名称"
$q4.Cells.Item(1, 4).Value = "基金规模"
$q4.Cells.Item(1, 5).Value = "股票总仓位"
$q4.Cells.Item(1, 6).Value = "仓位占比"
$q4.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q4.Cells.Item(1, 8).Value = "仓位排名"

# Column A (row index) on the data rows uses the same header style.
$summary.Range("A2").Copy($q4.Range("A2:A6")) | Out-Null

$rows = @(
    @(0, "161838", "银华创业板两年定期开放混合", "4.71", "95.23", "8.71", "0.4102", 4),
    @(1, "180020", "银华成长先锋混合",           "2.18", "78.61", "5.82", "0.1269", 5),
    @(2, "162203", "泰达宏利稳定混合",           "2.98", "91.72", "3.47", "0.1034", 9),
    @(3, "080001", "长盛成长价值混合A",          "2.43", "63.89", "3.36", "0.0816", 1),
    @(4, "012715", "长盛成长价值混合C",          "0.10", "63.89", "3.36", "0.0034", 1)
)

$r = 2
foreach ($row in $rows) {
    $q4.Cells.Item($r, 1).Value = $row[0]
    Set-TextCell $q4.Cells.Item($r, 2) $row[1]
    Set-TextCell $q4.Cells.Item($r, 3) $row[2]
    Set-TextCell $q4.Cells.Item($r, 4) $row[3]
    Set-TextCell $q4.Cells.Item($r, 5) $row[4]
    Set-TextCell $q4.Cells.Item($r, 6) $row[5]
    Set-TextCell $q4.Cells.Item($r, 7) $row[6]
    $q4.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# Leave the original "总计" sheet focused/active, matching the source
# workbook's view state.
$summary.Activate()
